# "Added ILSpy to resources"
#
# Slide 10 ("Resources") has a bulleted list of links. Add a new
# "ILSpy:" label line plus an indented hyperlinked URL line right after
# the existing "LINQ to Twitter" entry, keeping the trailing blank
# paragraph that already exists at the end of the placeholder.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$shp = $s.Shapes.Item("Content Placeholder 2")
$tf = $shp.TextFrame
$tr = $tf.TextRange

# "ILSpy" and ":" as two separate runs (two InsertAfter calls keep them
# as distinct runs instead of merging into one).
[void]$tr.InsertAfter("ILSpy")
[void]$tr.InsertAfter(":")
[void]$tr.InsertAfter("`r")

$urlText = "https://github.com/icsharpcode/ILSpy"
[void]$tr.InsertAfter($urlText)
[void]$tr.InsertAfter(" ")

# The URL line is a second-level bullet, like the other links on this slide.
$lastPara = $tr.Paragraphs($tr.Paragraphs().Count)
$lastPara.IndentLevel = 2

# Hyperlink just the URL text (not the trailing space).
$urlRange = $tr.Characters($lastPara.Start, $urlText.Length)
$urlRange.ActionSettings.Item(1).Hyperlink.Address = $urlText
